$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I9").Value = "sd"
$ws.Range("J9").Value = "Statement-non-opinion"
$ws.Range("I21").Value = "sv"
$ws.Range("J21").Value = "Statement-opinion"
$ws.Range("I23").Value = "b"
$ws.Range("J23").Value = "Acknowledge (Backchannel)"
$ws.Range("I32").Value = "ba"
$ws.Range("J32").Value = "Appreciation"
$ws.Range("I33").Value = "ba"
$ws.Range("J33").Value = "Appreciation"
$ws.Range("I36").Value = "b"
$ws.Range("J36").Value = "Acknowledge (Backchannel)"
$ws.Range("I40").Value = "sv"
$ws.Range("J40").Value = "Statement-opinion"
$ws.Range("I45").Value = "sv"
$ws.Range("J45").Value = "Statement-opinion"
$ws.Range("I47").Value = "aa"
$ws.Range("J47").Value = "Agree/Accept"
$ws.Range("I50").Value = "sv"
$ws.Range("J50").Value = "Statement-opinion"
$ws.Range("I58").Value = "sv"
$ws.Range("J58").Value = "Statement-opinion"
$ws.Range("I66").Value = "sv"
$ws.Range("J66").Value = "Statement-opinion"
$ws.Range("I79").Value = "sv"
$ws.Range("J79").Value = "Statement-opinion"
$ws.Range("I81").Value = "ba"
$ws.Range("J81").Value = "Appreciation"
$ws.Range("I82").Value = "sv"
$ws.Range("J82").Value = "Statement-opinion"
$ws.Range("I94").Value = "sv"
$ws.Range("J94").Value = "Statement-opinion"
$ws.Range("I96").Value = "sv"
$ws.Range("J96").Value = "Statement-opinion"
$ws.Range("I97").Value = "sd"
$ws.Range("J97").Value = "Statement-non-opinion"
$ws.Range("I111").Value = "sv"
$ws.Range("J111").Value = "Statement-opinion"
$ws.Range("I116").Value = "b"
$ws.Range("J116").Value = "Acknowledge (Backchannel)"
$ws.Range("I120").Value = "sv"
$ws.Range("J120").Value = "Statement-opinion"
$ws.Range("I135").Value = "sv"
$ws.Range("J135").Value = "Statement-opinion"
$ws.Range("I137").Value = "sv"
$ws.Range("J137").Value = "Statement-opinion"
$ws.Range("I142").Value = "aa"
$ws.Range("J142").Value = "Agree/Accept"
$ws.Range("I198").Value = "sv"
$ws.Range("J198").Value = "Statement-opinion"
$ws.Range("I200").Value = "ba"
$ws.Range("J200").Value = "Appreciation"
$ws.Range("I234").Value = "sv"
$ws.Range("J234").Value = "Statement-opinion"
$ws.Range("I302").Value = "sd"
$ws.Range("J302").Value = "Statement-non-opinion"
$ws.Range("I311").Value = "sv"
$ws.Range("J311").Value = "Statement-opinion"
$ws.Range("I327").Value = "%"
$ws.Range("J327").Value = "Uninterpretable"
$ws.Range("I352").Value = "ba"
$ws.Range("J352").Value = "Appreciation"
$ws.Range("I363").Value = "sd"
$ws.Range("J363").Value = "Statement-non-opinion"
$ws.Range("I364").Value = "%"
$ws.Range("J364").Value = "Uninterpretable"
$ws.Range("I371").Value = "ba"
$ws.Range("J371").Value = "Appreciation"
$ws.Range("I378").Value = "sd"
$ws.Range("J378").Value = "Statement-non-opinion"
$ws.Range("I387").Value = "b"
$ws.Range("J387").Value = "Acknowledge (Backchannel)"
$ws.Range("I393").Value = "sv"
$ws.Range("J393").Value = "Statement-opinion"
$ws.Range("I402").Value = "sv"
$ws.Range("J402").Value = "Statement-opinion"
$ws.Range("I404").Value = "sd"
$ws.Range("J404").Value = "Statement-non-opinion"
$ws.Range("I407").Value = "sv"
$ws.Range("J407").Value = "Statement-opinion"
$ws.Range("I422").Value = "sd"
$ws.Range("J422").Value = "Statement-non-opinion"
$ws.Range("I426").Value = "aa"
$ws.Range("J426").Value = "Agree/Accept"
$ws.Range("I429").Value = "sd"
$ws.Range("J429").Value = "Statement-non-opinion"
$ws.Range("I431").Value = "sd"
$ws.Range("J431").Value = "Statement-non-opinion"
